# Hospital_Records.xlsx edit script
# - renames sheet, trims Appointment Date / Visit Reason columns and the
#   last row, updates several patient rows, resizes the surviving columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab (space -> underscore).
$ws.Name = "Hospital_Records"

# 2. Drop the "Appointment Date" (F) and "Visit Reason" (G) columns entirely.
$ws.Columns("F:G").Delete()

# 3. Drop the last record (row 6 - Vikram Patil).
$ws.Rows("6").Delete()

# 4. Resize the remaining data columns. Excel's ColumnWidth property adds the
#    ~0.8333 (5/6) default-font padding on top of whatever is assigned, so we
#    back that constant out to land exactly on the requested stored widths.
$padding = 0.8333333333333334
$ws.Columns("A").ColumnWidth = 19 - $padding
$ws.Columns("B").ColumnWidth = 15 - $padding
$ws.Columns("C").ColumnWidth = 23 - $padding

# 5. Update the patient rows.
# Row 2 - Aarav Mehta keeps his name/address/age/gender, only the phone changes.
$ws.Range("B2").Value = "'+917823844614"

# Row 3 - now Vanshika panjwani.
$ws.Range("A3").Value = "Vanshika panjwani"
$ws.Range("B3").Value = "'+917823844614"
$ws.Range("C3").Value = "24 MG Road, Bengaluru"
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = "Male"

# Row 4 - now Vanshika panjwani.
$ws.Range("A4").Value = "Vanshika panjwani"
$ws.Range("B4").Value = "'+917823844614"
$ws.Range("C4").Value = "24 MG Road, Bengaluru"
$ws.Range("D4").Value = 28

# Row 5 - now Vanshika panjwani.
$ws.Range("A5").Value = "Vanshika panjwani"
$ws.Range("B5").Value = "'+918767545559"
$ws.Range("C5").Value = "24 MG Road, Bengaluru"
$ws.Range("D5").Value = 28
$ws.Range("E5").Value = "Male"

# The apostrophe prefixes above force the numeric-looking phone numbers to be
# stored as text (matching the source data) while adding a "quote prefix"
# style to those cells; strip that back off so the cells stay styled exactly
# like the rest of the sheet (no explicit style index).
$ws.Range("B2:B5").ClearFormats()
